$wb = $excel.ActiveWorkbook

# --- Sheet 1: ROW50-FE-LIFTER  (A1:I72 -> A1:I73), new row 73 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A73").Value = 45760.75132583333
$ws1.Range("A73").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B73").Value = "0x01,0x90"
$ws1.Range("C73").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Range("D73").Value = "0x01,0x4e"
$ws1.Range("E73").Value = "0xe"
$ws1.Range("F73").Value = 400
$ws1.Range("G73").Value = 568631262647114000000000.0
$ws1.Range("H73").Value = 334
$ws1.Range("I73").Value = 14

# --- Sheet 2: ROW50-MID-LIFTER  (A1:I74 -> A1:I75), new row 75 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A75").Value = 45760.71534722222
$ws2.Range("A75").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B75").Value = "0x01,0x90 "
$ws2.Range("C75").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Range("D75").Value = "0x01,0x56"
$ws2.Range("E75").Value = "0x19"
$ws2.Range("F75").Value = 400
$ws2.Range("G75").NumberFormat = "@"
$ws2.Range("G75").Value = "568631262647113771663628"
$ws2.Range("G75").ClearFormats()
$ws2.Range("H75").Value = 342
$ws2.Range("I75").Value = 25

# --- Sheet 3: ROW11-FE-LIFTER  (A1:I72 -> A1:I73), new row 73 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A73").Value = 45760.78471831018
$ws3.Range("A73").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("B73").Value = "0x01,0x90"
$ws3.Range("C73").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Range("D73").Value = "0x01,0x4e"
$ws3.Range("E73").Value = "0x14"
$ws3.Range("F73").Value = 400
$ws3.Range("G73").Value = 568631262647114000000000.0
$ws3.Range("H73").Value = 334
$ws3.Range("I73").Value = 20

# --- Sheet 4: ROW11-MID-LIFTER  (A1:I72 -> A1:I73), new row 73 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A73").Value = 45760.91061049768
$ws4.Range("A73").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Range("B73").Value = "0x01,0x90"
$ws4.Range("C73").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Range("D73").Value = "0x01,0x56"
$ws4.Range("E73").Value = "0x19"
$ws4.Range("F73").Value = 400
$ws4.Range("G73").Value = 568631262647114000000000.0
$ws4.Range("H73").Value = 342
$ws4.Range("I73").Value = 25
